# Updates cryptos list values (Price / Volume(1h) columns) per upstream data refresh.
# Rows 31-32 (Aptos / SuiNetwork) also swapped position, so their Coin/Link/Price/Volume all change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.484.93"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "2.275.81"
$ws.Range("E3").Value = "  -6.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'544.39"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'130.50"
$ws.Range("E6").Value = "  -5.22%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").Value = "'5.49"
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = "  -6.11%  "
$ws.Range("D13").Value = "'23.63"
$ws.Range("E13").Value = "  -6.05%  "
$ws.Range("D14").Value = "2.683.89"
$ws.Range("E14").Value = "  -6.09%  "
$ws.Range("D15").Value = "58.455.39"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "'0.0000133"
$ws.Range("E16").Value = "  -4.13%  "
$ws.Range("D17").Value = "2.272.19"
$ws.Range("E17").Value = "  -5.18%  "
$ws.Range("D18").Value = "'10.67"
$ws.Range("E18").Value = "  -5.89%  "
$ws.Range("D19").Value = "'4.29"
$ws.Range("E19").Value = "  -5.06%  "
$ws.Range("D20").Value = "'313.98"
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'62.86"
$ws.Range("E23").Value = "  -4.40%  "
$ws.Range("D24").Value = "'0.170"
$ws.Range("E24").Value = "  -5.01%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'8.09"
$ws.Range("E26").Value = "  -7.36%  "
$ws.Range("E27").Value = "  -6.21%  "
$ws.Range("D28").Value = "'1.74"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'170.99"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "0.0₃0719"
$ws.Range("E30").Value = "  -7.45%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").Value = "'1.07"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'5.76"
$ws.Range("E32").Value = "  -6.07%  "
$ws.Range("D33").Value = "'0.382"
$ws.Range("E33").Value = "  -6.22%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'17.72"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -6.21%  "
$ws.Range("D38").Value = "'3.95"
$ws.Range("E38").Value = "  -6.59%  "
$ws.Range("D39").Value = "'38.10"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("E40").Value = "  -6.24%  "
$ws.Range("D41").Value = "'299.32"
$ws.Range("E41").Value = "  -10.28%  "
$ws.Range("D42").Value = "'140.65"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = "  -6.11%  "
$ws.Range("D44").Value = "'0.0945"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "'0.0497"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("D46").Value = "'0.549"
$ws.Range("E46").Value = "  -5.26%  "
$ws.Range("D47").Value = "'18.28"
$ws.Range("E47").Value = "  -9.36%  "
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'16.46"
$ws.Range("E50").Value = "  -7.05%  "
$ws.Range("E51").Value = "  -0.47%  "
